$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.225.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.11%  "

$ws.Range("D3").Value = "'3.076.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.38%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'574.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.57%  "

$ws.Range("D6").Value = "'170.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.85%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "'3.076.18"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.21%  "

$ws.Range("D9").Value = "'0.510"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.21%  "

$ws.Range("D10").Value = "'6.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.87%  "

$ws.Range("E11").Value = "  -2.86%  "

$ws.Range("E12").Value = "  -2.57%  "

$ws.Range("E13").Value = "  -3.58%  "

$ws.Range("D14").Value = "'35.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.48%  "

$ws.Range("E15").Value = "  -1.44%  "

$ws.Range("D16").Value = "'3.587.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.44%  "

$ws.Range("D17").Value = "'66.148.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.18%  "

$ws.Range("E18").Value = "  -2.93%  "

$ws.Range("D19").Value = "'3.075.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.43%  "

$ws.Range("D20").Value = "'16.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.49%  "

$ws.Range("D21").Value = "'484.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.67%  "

$ws.Range("E22").Value = "  -3.39%  "

$ws.Range("E23").Value = "  -3.35%  "

$ws.Range("D24").Value = "'82.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.83%  "

$ws.Range("D25").Value = "'12.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.10%  "

$ws.Range("D26").Value = "'2.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.57%  "

$ws.Range("D27").Value = "'10.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.67%  "

$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("D29").Value = "'7.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.53%  "

$ws.Range("E30").Value = "  -5.40%  "

$ws.Range("E31").Value = "  -3.25%  "

$ws.Range("D32").Value = "'27.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.93%  "

$ws.Range("E33").Value = "  -4.20%  "

$ws.Range("D34").Value = "'0.0₃0900"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.89%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("B36").Value = "Arweave"
$ws.Range("C36").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D36").Value = "'47.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.34%  "

$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").Value = "'0.944"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.21%  "

$ws.Range("D38").Value = "'5.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.89%  "

$ws.Range("E39").Value = "  -1.00%  "

$ws.Range("E40").Value = "  -5.16%  "

$ws.Range("E41").Value = "  -4.08%  "

$ws.Range("D42").Value = "'8.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.97%  "

$ws.Range("D43").Value = "'2.776.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.44%  "

$ws.Range("D44").Value = "'2.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.50%  "

$ws.Range("E45").Value = "  -3.36%  "

$ws.Range("D46").Value = "'134.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.66%  "

$ws.Range("D47").Value = "'364.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.96%  "

$ws.Range("E48").Value = "  +0.02%  "

$ws.Range("D49").Value = "'24.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.91%  "

$ws.Range("D50").Value = "'2.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.48%  "

$ws.Range("E51").Value = "  -2.24%  "
